$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 112183921
$ws.Range("B2").Value = 89553
$ws.Range("Q2").Value = 763392
$ws.Range("R2").Value = 7448819

# Row 3
$ws.Range("A3").Value = 112183920
$ws.Range("B3").Value = 89553
$ws.Range("E3").Value = 1202
$ws.Range("F3").Value = "Ullticka"
$ws.Range("G3").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H3").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q3").Value = 763654
$ws.Range("R3").Value = 7448906

# Row 4
$ws.Range("A4").Value = 112181983
$ws.Range("B4").Value = 89718
$ws.Range("E4").Value = 1588
$ws.Range("F4").Value = "Violmussling"
$ws.Range("G4").Value = "Trichaptum laricinum"
$ws.Range("H4").Value = "(P.Karst.) Ryvarden"
$ws.Range("Q4").Value = 763400
$ws.Range("R4").Value = 7448829
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2023-07-06"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2023-07-06"

# Row 5
$ws.Range("A5").Value = 112183947
$ws.Range("B5").Value = 89718
$ws.Range("E5").Value = 1588
$ws.Range("F5").Value = "Violmussling"
$ws.Range("G5").Value = "Trichaptum laricinum"
$ws.Range("H5").Value = "(P.Karst.) Ryvarden"
$ws.Range("Q5").Value = 763391
$ws.Range("R5").Value = 7448820

# Row 6
$ws.Range("A6").Value = 112183036
$ws.Range("B6").Value = 89571
$ws.Range("E6").Value = 5432
$ws.Range("F6").Value = "Granticka"
$ws.Range("G6").Value = "Porodaedalea chrysoloma"
$ws.Range("H6").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q6").Value = 763482
$ws.Range("R6").Value = 7448939
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2023-06-30"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2023-06-30"

# Row 7
$ws.Range("A7").Value = 112181997
$ws.Range("B7").Value = 89553
$ws.Range("E7").Value = 1202
$ws.Range("F7").Value = "Ullticka"
$ws.Range("G7").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H7").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q7").Value = 763401
$ws.Range("R7").Value = 7448827
